# Added SCM log data
# Appends rows 34-42 to Sheet1, mirroring the new task-log entries that
# were recorded for the "future SCM" log (commit "Added SCM log data").
#
# The text cells are written in the same order the original author typed
# them in (columns A/B filled in first for a few rows, then circling back
# to fill in column C for earlier rows), so that the shared-string table
# ends up in the same order as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34 ---
$ws.Range("A34").Value = "Added test case reference"
$ws.Range("B34").Value = "CID e5f342d"

# --- Row 35 ---
$ws.Range("A35").Value = "Added cell # and assn 10"
$ws.Range("B35").Value = "CID 5d2f5d9"
$ws.Range("C35").Value = "Phone # and doc"

# --- Row 36 ---
$ws.Range("A36").Value = "Updated assn 10"
$ws.Range("B36").Value = "CID 8c3aed0"
$ws.Range("C36").Value = "document"

# --- back to Row 34 ---
$ws.Range("C34").Value = "Tracebility matrix"

# --- Row 37 ---
$ws.Range("A37").Value = "Added test case design "
$ws.Range("B37").Value = "CID 5s70281"

# --- Row 38 ---
$ws.Range("A38").Value = "Revised test case design"
$ws.Range("B38").Value = "CID d7aa179"
$ws.Range("C38").Value = "Document revision"

# --- back to Row 37 ---
$ws.Range("C37").Value = "Test case design.xls"

# --- Row 39 ---
$ws.Range("A39").Value = "Proposed change mangement diagram"
$ws.Range("B39").Value = "CID 8d4e101"
$ws.Range("C39").Value = "Change mangement.ppt"

# --- Row 40 ---
$ws.Range("A40").Value = "Proposed assignment 10 .ppt"
$ws.Range("B40").Value = "CID 317401e"
$ws.Range("C40").Value = "Assn 10 .ppt"

# --- Row 41 ---
$ws.Range("A41").Value = "Updated assn 10 and traceability matrix"
$ws.Range("C41").Value = "Revised existing docs"
$ws.Range("B41").Value = "CID 237a8c4"
$ws.Range("G41").Value = "In care of Susan"

# --- Row 42 ---
$ws.Range("A42").Value = "Created RSC spreadsheet for GP Genie"
$ws.Range("B42").Value = "CID 870b9d0"
$ws.Range("C42").Value = "SW measurement tool"

# --- Where (date), When (who), How Long (hours) columns ---
$dates = @{ 34=41610; 35=41610; 36=41611; 37=41611; 38=41613; 39=41615; 40=41615; 41=41615; 42=41616 }
$hours = @{ 34=0.5;   35=2;     36=0.25;  37=1;     38=0.75;  39=0.75;  40=1;     41=2;     42=2 }

foreach ($rowNum in 34..42) {
    $ws.Range("D$rowNum").Value = $dates[$rowNum]

    # Reuse the existing m/d/yyyy number format (style index 9) used by the
    # rest of the "Where" column instead of creating a new numFmt entry.
    $ws.Range("D26").Copy()
    $ws.Range("D$rowNum").PasteSpecial(-4122)

    $ws.Range("E$rowNum").Value = "Roger"
    $ws.Range("F$rowNum").Value = $hours[$rowNum]
}

$excel.CutCopyMode = 0

$ws.Range("E43").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
